# Freelance Timesheet - add a new session entry (row 17) describing the
# "session history / current availability" query support, and update the
# selection to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timesheet row: 2023-02-20, 10:08am - 11:08am, 1 hour.
$ws.Range("A17").NumberFormat = "d-mmm"
$ws.Range("A17").Value = 44977
$ws.Range("B17").NumberFormat = "h:mm"
$ws.Range("B17").Value = 0.42222222222222222
$ws.Range("C17").NumberFormat = "h:mm"
$ws.Range("C17").Value = 0.46388888888888885
$ws.Range("D17").Value = "Support for simple queries: session history for both students and tutors, current weekly availability for both students and tutors."
$ws.Range("E17").Value = 1

# Row height grows to fit the wrapped description text (matches the other
# multi-line rows in the log).
$ws.Rows(17).RowHeight = 56

# Cursor ends up on A18 after entering the row.
$ws.Range("A18").Select()
